$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "Duke 18-19" worksheet after the last existing sheet ("Gale 18-19")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Duke 18-19"

# ---------------------------------------------------------------------------
# Column widths (approximate character-width equivalents of the source
# report's raw widths: A:B=36, C:J=10.71, K:V=5.43)
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 36
$ws.Range("C1:J1").EntireColumn.ColumnWidth = 9.9
$ws.Range("K1:V1").EntireColumn.ColumnWidth = 4.6

# ---------------------------------------------------------------------------
# Title / report metadata block (rows 1-7)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Journal Report 1 (R4)"
$ws.Range("B1").Value = "Number of Successful Full-Text Article Requests by Month and Journal"
$ws.Range("A2").Value = "FLORIDA STATE UNIV"
$ws.Range("A3").Value = " "
$ws.Range("A4").Value = "Period covered by Report:"
$ws.Range("A5").Value = "2018-07-01 to 2019-06-30"
$ws.Range("A6").Value = "Date run:"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2019-09-05"

$titleBlock = $ws.Range("A1:V7")
$titleBlock.Font.Name = "Arial"
$titleBlock.Font.Size = 9
$titleBlock.Font.Bold = $true
$titleBlock.Font.ColorIndex = 1
$titleBlock.HorizontalAlignment = -4131
$titleBlock.VerticalAlignment = -4160

$ws.Rows("1").RowHeight = 15.95
$ws.Rows("2").RowHeight = 15.95
$ws.Rows("3").RowHeight = 15.95
$ws.Rows("4").RowHeight = 15.95
$ws.Rows("5").RowHeight = 15.95
$ws.Rows("6").RowHeight = 15.95
$ws.Rows("7").RowHeight = 12

# ---------------------------------------------------------------------------
# Column header row (row 8)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Journal"
$ws.Range("B8").Value = "Publisher"
$ws.Range("C8").Value = "Platform"
$ws.Range("D8").Value = "Journal DOI"
$ws.Range("E8").Value = "Proprietary Identifier"
$ws.Range("F8").Value = "Print ISSN"
$ws.Range("G8").Value = "Online ISSN"
$ws.Range("H8").Value = "Reporting Period Total"
$ws.Range("I8").Value = "Reporting Period HTML"
$ws.Range("J8").Value = "Reporting Period PDF"
$ws.Range("K8").Value = "Jul-2018"
$ws.Range("L8").Value = "Aug-2018"
$ws.Range("M8").Value = "Sep-2018"
$ws.Range("N8").Value = "Oct-2018"
$ws.Range("O8").Value = "Nov-2018"
$ws.Range("P8").Value = "Dec-2018"
$ws.Range("Q8").Value = "Jan-2019"
$ws.Range("R8").Value = "Feb-2019"
$ws.Range("S8").Value = "Mar-2019"
$ws.Range("T8").Value = "Apr-2019"
$ws.Range("U8").Value = "May-2019"
$ws.Range("V8").Value = "Jun-2019"

$headerRow = $ws.Range("A8:V8")
$headerRow.Font.Name = "Arial"
$headerRow.Font.Size = 9
$headerRow.Font.Bold = $true
$headerRow.Font.ColorIndex = 1
$headerRow.Interior.Color = 11439959
$headerRow.HorizontalAlignment = -4131

$ws.Range("E8").WrapText = $true
$ws.Range("K8:V8").WrapText = $true

$ws.Rows("8").RowHeight = 45

# ---------------------------------------------------------------------------
# Total row (row 9)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Total for all journals"
$ws.Range("C9").Value = "Silverchair"

$totalRow = $ws.Range("A9:V9")
$totalRow.Font.Name = "Arial"
$totalRow.Font.Size = 9
$totalRow.Font.ColorIndex = 1
$totalRow.Interior.Color = 14606524
$totalRow.HorizontalAlignment = -4131
$totalRow.VerticalAlignment = -4160
$totalRow.WrapText = $true

$ws.Rows("9").RowHeight = 15

# ---------------------------------------------------------------------------
# Data rows (rows 10-14)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Camera Obscura: Feminism, Culture, and Media Studies"
$ws.Range("B10").Value = "Duke University Press"
$ws.Range("C10").Value = "Silverchair"
$ws.Range("F10").Value = "0270-5346"
$ws.Range("G10").Value = "1529-1510"
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 2
$ws.Cells.Item(10, 15).Value = 6
$ws.Cells.Item(10, 16).Value = 8
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 6
$ws.Cells.Item(10, 19).Value = 6
$ws.Cells.Item(10, 20).Value = 4
$ws.Cells.Item(10, 21).Value = 0
$ws.Cells.Item(10, 22).Value = 0

$ws.Range("A11").Value = "Environmental Humanities"
$ws.Range("B11").Value = "Duke University Press"
$ws.Range("C11").Value = "Silverchair"
$ws.Range("F11").Value = "2201-1919"
$ws.Range("G11").Value = "2201-1919"
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 2
$ws.Cells.Item(11, 14).Value = 4
$ws.Cells.Item(11, 15).Value = 2
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 8
$ws.Cells.Item(11, 20).Value = 4
$ws.Cells.Item(11, 21).Value = 0
$ws.Cells.Item(11, 22).Value = 4

$ws.Range("A12").Value = "Journal of Health Politics, Policy and Law"
$ws.Range("B12").Value = "Duke University Press"
$ws.Range("C12").Value = "Silverchair"
$ws.Range("F12").Value = "0361-6878"
$ws.Range("G12").Value = "1527-1927"
$ws.Cells.Item(12, 11).Value = 6
$ws.Cells.Item(12, 12).Value = 6
$ws.Cells.Item(12, 13).Value = 36
$ws.Cells.Item(12, 14).Value = 16
$ws.Cells.Item(12, 15).Value = 8
$ws.Cells.Item(12, 16).Value = 14
$ws.Cells.Item(12, 17).Value = 6
$ws.Cells.Item(12, 18).Value = 18
$ws.Cells.Item(12, 19).Value = 16
$ws.Cells.Item(12, 20).Value = 18
$ws.Cells.Item(12, 21).Value = 14
$ws.Cells.Item(12, 22).Value = 0

$ws.Range("A13").Value = "Journal of Korean Studies"
$ws.Range("B13").Value = "Duke University Press"
$ws.Range("C13").Value = "Silverchair"
$ws.Range("F13").Value = "2158-1665"
$ws.Range("G13").Value = "0731-1613"
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 0
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(13, 21).Value = 0
$ws.Cells.Item(13, 22).Value = 0

$ws.Range("A14").Value = "the minnesota review"
$ws.Range("B14").Value = "Duke University Press"
$ws.Range("C14").Value = "Silverchair"
$ws.Range("F14").Value = "0026-5667"
$ws.Range("G14").Value = "2157-4189"
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 4
$ws.Cells.Item(14, 15).Value = 6
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 0
$ws.Cells.Item(14, 21).Value = 0
$ws.Cells.Item(14, 22).Value = 0

$dataRows = $ws.Range("A10:V14")
$dataRows.Font.Name = "Arial"
$dataRows.Font.Size = 9
$dataRows.Font.ColorIndex = 1
$dataRows.HorizontalAlignment = -4131
$dataRows.VerticalAlignment = -4160
$dataRows.WrapText = $true

$ws.Rows("10").RowHeight = 15
$ws.Rows("11").RowHeight = 15
$ws.Rows("12").RowHeight = 15
$ws.Rows("13").RowHeight = 15
$ws.Rows("14").RowHeight = 15

# ---------------------------------------------------------------------------
# Page setup / print options, matching the source report
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintGridlines = $true
$ws.PageSetup.Orientation = 2

# ---------------------------------------------------------------------------
# Selection / view state - the new sheet becomes the active (selected) tab
# ---------------------------------------------------------------------------
$ws.Range("H34").Select() | Out-Null
